# "handling login functionality using excelsheet"
#
# This script reproduces, via Excel COM-interop calls, the authored edit:
#   - rename Sheet2 -> "ali" and make it the active/visible tab
#   - select C3 on Sheet1 (no longer the active tab) and D4 on "ali"
#   - add four new rows of data (with shared strings) to "ali"
#   - size the two columns on "ali"
#   - style the new "jobprogram" cell (Arial 8, dark grey, centered + wrapped)
#   - set page setup (paper size / orientation) on "ali"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- rename Sheet2 ---------------------------------------------------------
$ws2.Name = "ali"

# --- column widths on the renamed sheet ------------------------------------
# (COM ColumnWidth is quantized by this host to ~1/6-character steps, so
# these inputs are chosen to land on the grid step nearest the authored
# widths of 32.44140625 / 17.6640625 characters)
$ws2.Columns.Item(1).ColumnWidth = 31.6
$ws2.Columns.Item(2).ColumnWidth = 16.77

# --- new rows of data -------------------------------------------------------
$ws2.Range("A3").Value = "abc"
$ws2.Range("B3").Value = "fks"

$ws2.Range("A4").Value = "dsd"
$ws2.Range("B4").Value = "fcsdc"

$ws2.Range("A5").Value = "training@jalaacademy.com"
$ws2.Range("B5").Value = "jobprogram"

$ws2.Range("A6").Value = "fghgh"
$ws2.Range("B6").Value = "dfggd"

# --- style the B5 "jobprogram" cell -----------------------------------------
$b5 = $ws2.Range("B5")
$b5.Font.Name = "Arial"
$b5.Font.Size = 8
$b5.Font.Color = 3355443   # RGB(51,51,51) -> 0x333333
$b5.HorizontalAlignment = -4108  # xlCenter
$b5.VerticalAlignment = -4108    # xlCenter
$b5.WrapText = $true

# --- page setup for "ali" ---------------------------------------------------
$ws2.PageSetup.PaperSize = 9       # xlPaperA4
$ws2.PageSetup.Orientation = 1     # xlPortrait

# --- selections / active tab ------------------------------------------------
[void]$ws1.Select()
[void]$ws1.Range("C3").Select()

[void]$ws2.Select()
[void]$ws2.Range("D4").Select()
